$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete column A entirely - this shifts B:F left to A:E, which removes
# the old un-headered first column (0/18 values) and drops the old
# trailing GENE column's duplicate position, matching the new A1:E3 layout.
$ws.Range("A1").EntireColumn.Delete()
